# jeden dodany wiersz bbb
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

# Turn on iterative calculation (calcPr calcId="145621" iterateCount="1")
$excel.Iteration = $true
$excel.MaxIterations = 1

# Add the new row of data at the bottom of the table (A18:B18)
$ws.Range("A18").Value = "bbb"
$ws.Range("B18").Value = "bbbb"

# Move/restore the active selection to B17
$ws.Range("B17").Select()
